$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.500607371330261
$ws.Range("B1").Value = 2.539062738418579
$ws.Range("C1").Value = 3.026667356491089
$ws.Range("D1").Value = 3.336366891860962
$ws.Range("E1").Value = 1.181021451950073
